$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.528.21"
$ws.Range("E2").Value = "  +3.31%  "

$ws.Range("D3").Value = "3.347.69"
$ws.Range("E3").Value = "  +7.75%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'246.35"
$ws.Range("E5").Value = "  +3.42%  "

$ws.Range("D6").Value = "'624.35"
$ws.Range("E6").Value = "  +1.10%  "

$ws.Range("E7").Value = "  +0.67%  "

$ws.Range("D8").Value = "'0.389"
$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "3.350.06"
$ws.Range("E10").Value = "  +7.80%  "

$ws.Range("D11").Value = "'0.801"
$ws.Range("E11").Value = "  -3.68%  "

$ws.Range("D12").Value = "'0.201"
$ws.Range("E12").Value = "  +1.59%  "

$ws.Range("D13").Value = "97.397.55"
$ws.Range("E13").Value = "  +3.84%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'35.83"
$ws.Range("E14").Value = "  +2.03%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000249"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").Value = "3.968.34"
$ws.Range("E16").Value = "  +7.81%  "

$ws.Range("D17").Value = "'5.55"
$ws.Range("E17").Value = "  +2.43%  "

$ws.Range("D18").Value = "3.360.53"
$ws.Range("E18").Value = "  +8.45%  "

$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").Value = "'15.30"
$ws.Range("E20").Value = "  +2.42%  "

$ws.Range("D21").Value = "'493.36"
$ws.Range("E21").Value = "  +10.46%  "

$ws.Range("D22").Value = "'0.0000212"
$ws.Range("E22").Value = "  +6.02%  "

$ws.Range("D23").Value = "'5.93"
$ws.Range("E23").Value = "  -1.05%  "

$ws.Range("D24").Value = "'9.33"
$ws.Range("E24").Value = "  +3.63%  "

$ws.Range("D25").Value = "'5.71"
$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("D26").Value = "'88.80"
$ws.Range("E26").Value = "  +3.11%  "

$ws.Range("D27").Value = "'12.19"
$ws.Range("E27").Value = "  -0.63%  "

$ws.Range("D28").Value = "3.504.17"
$ws.Range("E28").Value = "  +7.16%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("D30").Value = "'0.183"
$ws.Range("E30").Value = "  +1.66%  "

$ws.Range("E31").Value = "  -1.06%  "

$ws.Range("E32").Value = "  -1.79%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").Value = "'9.42"
$ws.Range("E34").Value = "  +1.28%  "

$ws.Range("D35").Value = "'27.83"
$ws.Range("E35").Value = "  +6.53%  "

$ws.Range("E36").Value = "  -4.78%  "

$ws.Range("D37").Value = "'7.56"
$ws.Range("E37").Value = "  -4.56%  "

$ws.Range("D38").Value = "'505.04"
$ws.Range("E38").Value = "  +5.59%  "

$ws.Range("E39").Value = "  +2.31%  "

$ws.Range("D40").Value = "'24.72"
$ws.Range("E40").Value = "  +3.03%  "

$ws.Range("D41").Value = "'0.454"
$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").Value = "'1.29"
$ws.Range("E42").Value = "  +0.30%  "

$ws.Range("D43").Value = "'3.31"
$ws.Range("E43").Value = "  +2.01%  "

$ws.Range("D44").Value = "'0.804"
$ws.Range("E44").Value = "  +16.44%  "

$ws.Range("E45").Value = "  -7.06%  "

$ws.Range("D47").Value = "'160.96"
$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("E48").Value = "  +5.08%  "

$ws.Range("D49").Value = "'4.61"
$ws.Range("E49").Value = "  +3.74%  "

$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "'1.36"
$ws.Range("E50").Value = "  +3.14%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0334"
$ws.Range("E51").Value = "  +3.39%  "
